$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty cells F3 and F24 (they had no value, just an empty
# inlineStr cell in the original markup).
$ws.Range("F3").ClearContents()
$ws.Range("F24").ClearContents()

# Correct the district name "Kalaburgi" -> "Kalaburagi (Gulbarga)" on every
# row where it appears verbatim in column G.
$rows = @(7, 8, 11, 13, 14, 15, 18, 20, 21, 22, 23, 25, 29, 30, 33, 34, 35, 36, 39)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Kalaburagi (Gulbarga)"
}
